$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (X) after the existing last data column (W),
# mirroring the formatting of column W, then fill in the 2020 values.
$ws.Range("W4:W16").Copy()
$ws.Range("X4").PasteSpecial(-4122)

$ws.Range("X4").Value = 2020
$ws.Range("X5").Value = 45.3
$ws.Range("X6").Value = 48.2
$ws.Range("X7").Value = 43.6
$ws.Range("X8").Value = 48.8
$ws.Range("X9").Value = 41.5
$ws.Range("X10").Value = 49.7
$ws.Range("X11").Value = 46.7
$ws.Range("X12").Value = 36.5
$ws.Range("X13").Value = 29.6
$ws.Range("X14").Value = 54.7
$ws.Range("X15").Value = 51.6
$ws.Range("X16").Value = 47.2

# Move the active selection the way it ended up in the saved workbook.
$ws.Range("AI21").Select() | Out-Null
